$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("M2").Value = 139
$ws.Range("M3").Value = 170
$ws.Range("M4").Value = 34
$ws.Range("M5").Value = 8
$ws.Range("M6").Value = 126
$ws.Range("M7").Value = 477

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("M2").Value = 2
$ws.Range("M7").Value = 4

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("M2").Value = 9
$ws.Range("M3").Value = 10
$ws.Range("M4").Value = 1
$ws.Range("M7").Value = 30

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("M3").Value = 8
$ws.Range("M6").Value = 8
$ws.Range("M7").Value = 20

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("M3").Value = 8
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 21

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("M7").Value = 16
$ws.Range("M8").Value = 30
$ws.Range("M9").Value = 5
$ws.Range("M11").Value = 5
$ws.Range("M15").Value = 5
$ws.Range("M19").Value = 14
$ws.Range("M20").Value = 19
$ws.Range("M23").Value = 7
$ws.Range("M24").Value = 4
$ws.Range("M28").Value = 25
$ws.Range("M30").Value = 8
$ws.Range("M32").Value = 20
$ws.Range("M36").Value = 21
$ws.Range("M41").Value = 12
$ws.Range("M46").Value = 4
$ws.Range("M47").Value = 11
$ws.Range("M52").Value = 4
$ws.Range("M53").Value = 12
$ws.Range("M60").Value = 2
$ws.Range("M61").Value = 1
$ws.Range("M65").Value = 15
$ws.Range("M71").Value = 5
$ws.Range("M75").Value = 6
$ws.Range("M76").Value = 8
$ws.Range("M77").Value = 13
$ws.Range("M81").Value = 12
$ws.Range("M83").Value = 25
$ws.Range("M86").Value = 4
$ws.Range("M88").Value = 8
$ws.Range("M89").Value = 5
$ws.Range("M91").Value = 5
$ws.Range("M99").Value = 477

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("M3").Value = 2
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("M2").Value = 4
$ws.Range("M3").Value = 7
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("M5").Value = 4
$ws.Range("M6").Value = 12

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("M2").Value = 8
$ws.Range("M3").Value = 7
$ws.Range("M4").Value = 4
$ws.Range("M7").Value = 25

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("M5").Value = 6
$ws.Range("M6").Value = 11

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("M3").Value = 8
$ws.Range("M6").Value = 4
$ws.Range("M7").Value = 14

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("M3").Value = 7
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("M3").Value = 4
$ws.Range("M6").Value = 8

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J3").Value = 3
$ws.Range("J5").Value = 4

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("M4").Value = 1
$ws.Range("M6").Value = 7

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("M4").Value = 1
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("M2").Value = 6
$ws.Range("M7").Value = 13

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("M2").Value = 8
$ws.Range("M5").Value = 1
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K2").Value = 3
$ws.Range("K6").Value = 5

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("M2").Value = 6
$ws.Range("M7").Value = 16

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("M3").Value = 2
$ws.Range("M7").Value = 4

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("M2").Value = 2
$ws.Range("M6").Value = 5

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("M3").Value = 1
$ws.Range("M7").Value = 5

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("M2").Value = 3
$ws.Range("M6").Value = 5

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 4

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 8

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("M3").Value = 12
$ws.Range("M5").Value = 2
$ws.Range("M7").Value = 25

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 4
$ws.Range("L6").Value = 6

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 2
